# Uniswap tracker update: add two new observation rows (27, 28) and two new
# derived-metric columns (M "delta vs fact", N "delta vs fact %") down the
# whole table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) New rows 27 & 28 - same shape as the existing row 26 (date/price/
#    plan/fact columns F, G, J, K). Copy row 26's number formats down
#    first so the new date cells pick up the short-date style, etc.
# ---------------------------------------------------------------------
$ws.Range("F26:K26").Copy() | Out-Null
$ws.Range("F27:K28").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("F27").Value = 45625      # 2024-11-29
$ws.Range("G27").Value = 3571.93
$ws.Range("J27").Value = 48
$ws.Range("K27").Value = 500.74

$ws.Range("F28").Value = 45630      # 2024-12-04
$ws.Range("G28").Value = 3669.19
$ws.Range("J28").Value = 60
$ws.Range("K28").Value = 634.55

# ---------------------------------------------------------------------
# 2) Column H ("verh. - D12" / premium vs floor) and column I (H / D14)
#    now need to reach down into the two new rows.
# ---------------------------------------------------------------------
$ws.Range("H12:H25").Formula = "=G12-`$D`$12"
$ws.Range("H26").Formula = "=G26-`$D`$12"
$ws.Range("H27").Formula = "=G27-`$D`$12"
$ws.Range("H28").Formula = "=G28-`$D`$12"

$ws.Range("I12:I28").Formula = "=H12/`$D`$14"

# ---------------------------------------------------------------------
# 3) New columns M (J-I, "факт vs расчет delta") and N ((J/I-1)*100,
#    "факт vs расчет delta %") for every data row 12..28.
# ---------------------------------------------------------------------
$ws.Range("M12").Formula = "=J12-I12"
$ws.Range("N12").Formula = "=(J12/I12-1)*100"

$ws.Range("M13:M28").Formula = "=J13-I13"
$ws.Range("N13:N28").Formula = "=(J13/I13-1)*100"

# ---------------------------------------------------------------------
# 4) Move the on-screen selection to where the author left off editing.
# ---------------------------------------------------------------------
$ws.Range("H28").Select() | Out-Null
